# Dataset variables explainer.xlsx — update the "AllGazeData" sheet:
#  - clarify the "outofbounds" explanation text
#  - remove the obsolete "smoothed_gaze" variable row
#  - leave the sheet active/selected (matches the saved workbook state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllGazeData")

# Row 5 ("outofbounds"): clarify that the gaze bounds refer to the screen
$ws.Range("B5").Value = "fraction of the samples where the gaze is out of bounds of the screen"

# Row 6 ("smoothed_gaze") is no longer documented — delete the whole row,
# shifting the remaining rows (i_leftreward, i_rightreward, i_leftcost x2) up.
$ws.Rows.Item(6).Delete()

# The workbook was last saved with "AllGazeData" as the active sheet/tab,
# with the final data row selected.
$ws.Activate()
$ws.Range("B10").Select()
